$d = $word.ActiveDocument

# --- 1) Justify the two header merge-field paragraphs ----------------------
# Paragraph containing the "header_position_tier_0" MERGEFIELD.
$r1 = $d.Content
$r1.Find.Execute("header_position_tier_0", $false, $false, $false, $false, $false, `
                  $true, 1, $false, "", 0)
if ($r1.Find.Found) {
    $r1.Paragraphs(1).Alignment = 3   # wdAlignParagraphJustify
}

# Paragraph containing the "header_rank_tier_0" / "header_username_tier_0"
# MERGEFIELDs.
$r2 = $d.Content
$r2.Find.Execute("header_rank_tier_0", $false, $false, $false, $false, $false, `
                  $true, 1, $false, "", 0)
if ($r2.Find.Found) {
    $r2.Paragraphs(1).Alignment = 3   # wdAlignParagraphJustify
}

# --- 2) Re-anchor the edit-position bookmark inside the footer_username ----
#        merge-field result. The user clicked/edited right after
#        "«footer_us" while clearing the widget's placeholder/help text, so
#        Word splits that text run in two around its automatic "_GoBack"
#        bookmark (and removes the old "_GoBack" bookmark that used to sit
#        at the very end of the document).
$r3 = $d.Content
$r3.Find.Execute("footer_username_tier_0", $false, $false, $false, $false, $false, `
                  $true, 1, $false, "", 0)
if ($r3.Find.Found) {
    $splitPoint = $d.Range($r3.Start + 9, $r3.Start + 9)
    $d.Bookmarks.Add("_GoBack", $splitPoint)
}
